$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values to reflect the reshuffled content
$ws.Range("B1").Value = 'Ementa atual:'
$ws.Range("C1").Value = 'Ementa modificada (dados modificados em vermelho):'
$ws.Range("B2").Value = 'LOQ4064'
$ws.Range("C2").Value = 'LOQ4064'
$ws.Range("A3").Value = 'Nome:'
$ws.Range("B3").Value = ' Engenharia de Processos Quimicos I'
$ws.Range("C3").Value = ' Engenharia de Processos Quimicos I'
$ws.Range("A4").Value = 'Name:'
$ws.Range("B4").Value = 'Chemical Process Engineering I'
$ws.Range("C4").Value = 'Chemical Process Engineering I'
$ws.Range("A5").Value = 'Créditos-aula:'
$ws.Range("B5").Value = '2'
$ws.Range("C5").Value = '2'
$ws.Range("A6").Value = 'Créditos-trabalho'
$ws.Range("B6").Value = '4'
$ws.Range("C6").Value = '4'
$ws.Range("A7").Value = 'Carga horária:'
$ws.Range("B7").Value = '150 h'
$ws.Range("C7").Value = '150 h'
$ws.Range("A8").Value = 'Ativação:'
$ws.Range("B8").Value = '01/01/2016'
$ws.Range("C8").Value = '01/01/2016'
$ws.Range("A9").Value = 'Semestre ideal:'
$ws.Range("B9").Value = 'EQD-7,EQN-8'
$ws.Range("C9").Value = 'EQD-7,EQN-8'
$ws.Range("A10").Value = 'Objetivos:'
$ws.Range("B10").Value = '5816812 - João Paulo Alves Silva'
$ws.Range("C10").Value = '5816812 - João Paulo Alves Silva'
$ws.Range("A11").Value = 'Objectives:'
$ws.Range("B11").Value = '1 - Consolidation and application of knowledge acquired in each of the specific areas of the Chemical Engineering degree. 
2 - Integration of knowledge of Chemical Engineering'
$ws.Range("C11").Value = '1 - Consolidation and application of knowledge acquired in each of the specific areas of the Chemical Engineering degree. 
2 - Integration of knowledge of Chemical Engineering'
$ws.Range("A12").Value = 'Docentes responsáveis:'
$ws.Range("A13").Value = 'Programa resumido:'
$ws.Range("B13").Value = 'Semestral'
$ws.Range("C13").Value = 'Semestral'
$ws.Range("A14").Value = 'Short syllabus:'
$ws.Range("B14").Value = '1 - Diagrams for Understanding Chemical Processes. 2 - The Structure and Synthesis of Process Flow Diagrams. 3 - Analysis of process performance. 4 - Industrial chemical plant study'
$ws.Range("C14").Value = '1 - Diagrams for Understanding Chemical Processes. 2 - The Structure and Synthesis of Process Flow Diagrams. 3 - Analysis of process performance. 4 - Industrial chemical plant study'
$ws.Range("A15").Value = 'Programa:'
$ws.Range("B15").Value = '01/01/2016'
$ws.Range("C15").Value = '01/01/2016'
$ws.Range("A16").Value = 'Syllabus:'
$ws.Range("B16").Value = '1 - Diagrams for Understanding Chemical Processes: Block Flow Diagrams; Process Flow Diagram (PFD); Piping and Instrumentation Diagram (P&ID).
2 - The Structure and Synthesis of Process Flow Diagrams:  Hierarchy of Process Design; Step 1 - Batch versus Continuous Process; Step 2 - The Input/Output Structure of the Process; Step 3 - The Recycle Structure of the Process
3 - Analysis of process performance: Process Input/Output Models; Tools for evaluating process performance.
4 - Industrial chemical plant study.'
$ws.Range("C16").Value = '1 - Diagrams for Understanding Chemical Processes: Block Flow Diagrams; Process Flow Diagram (PFD); Piping and Instrumentation Diagram (P&ID).
2 - The Structure and Synthesis of Process Flow Diagrams:  Hierarchy of Process Design; Step 1 - Batch versus Continuous Process; Step 2 - The Input/Output Structure of the Process; Step 3 - The Recycle Structure of the Process
3 - Analysis of process performance: Process Input/Output Models; Tools for evaluating process performance.
4 - Industrial chemical plant study.'
$ws.Range("A17").Value = 'Avaliação:'
$ws.Range("A18").Value = 'Método:'
$ws.Range("B18").Value = '5816812 - João Paulo Alves Silva'
$ws.Range("C18").Value = '5816812 - João Paulo Alves Silva'
$ws.Range("A19").Value = 'Critério:'
$ws.Range("B19").Value = 'Provas escritas e Apresentação de Trabalhos'
$ws.Range("C19").Value = 'Provas escritas e Apresentação de Trabalhos'
$ws.Range("A20").Value = 'Norma de recuperação:'
$ws.Range("B20").Value = 'A nota será composta por ao menos uma prova escrita e trabalhos realizados e apresentados durante o semestre. O peso de cada atividade será definido segundo critérios do professor.'
$ws.Range("C20").Value = 'A nota será composta por ao menos uma prova escrita e trabalhos realizados e apresentados durante o semestre. O peso de cada atividade será definido segundo critérios do professor.'
$ws.Range("A21").Value = 'Bibliografia:'
$ws.Range("B21").Value = 'Média Final = (N + Prova Recuperação)/2'
$ws.Range("C21").Value = 'Média Final = (N + Prova Recuperação)/2'
$ws.Range("A22").Value = 'Requisitos:'
$ws.Range("B23").Value = 'LOQ4054 -  Fenômenos de Transporte III  (Requisito fraco)
'
$ws.Range("C23").Value = 'LOQ4054 -  Fenômenos de Transporte III  (Requisito fraco)
'
$ws.Range("B24").Value = 'LOQ4086 -  Operações Unitárias II  (Requisito fraco)
'
$ws.Range("C24").Value = 'LOQ4086 -  Operações Unitárias II  (Requisito fraco)
'

# Clear cells that no longer hold content
$ws.Range("B17").ClearContents()
$ws.Range("C17").ClearContents()
$ws.Range("B22").ClearContents()
$ws.Range("C22").ClearContents()
$ws.Range("A23").ClearContents()

# Remove the now-unused trailing row (shifts dimension from C25 to C24)
$ws.Rows.Item(25).Delete()

# Row height adjustments for rows whose wrapped content changed size
$ws.Rows.Item(13).RowHeight = 60
$ws.Rows.Item(15).RowHeight = 120
$ws.Rows.Item(17).RowHeight = 15
$ws.Rows.Item(18).RowHeight = 60
$ws.Rows.Item(21).RowHeight = 120
$ws.Rows.Item(22).RowHeight = 15
$ws.Rows.Item(23).RowHeight = 30
